$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# Row 2 new formulas
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Row 3
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Rows 4-15 area formula, shared group
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# D4:D9 shared formula range shrink (previously D4:D12)
$ws.Range("D4:D9").Formula = "=(A4/100+(A5/100-A4/100)/2)"

# E6:E9 shared formula range shrink (previously E6:E12)
$ws.Range("E6:E9").Formula = "=(D6-D5)*(B6/100)*C6"

# Selection and view settings
$ws.Range("J2:K2").Select()

$wb.Save()
